$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "unab@unab.edu.co"
$ws.Range("B2").Value = "2023-05-18 15:03:40"
$ws.Range("D2").Value = 22222

$ws.Range("A3").Value = "unab@unab.edu.co"
$ws.Range("B3").Value = "2023-05-18 16:35:07"
$ws.Range("C3").Value = "Educacion"
$ws.Range("D3").Value = 150000

$ws.Range("A4").Value = "unab@unab.edu.co"
$ws.Range("B4").Value = "2023-05-18 16:35:17"
$ws.Range("C4").Value = "Alojamiento"
$ws.Range("D4").Value = 200000

$ws.Range("A5").Value = "unab@unab.edu.co"
$ws.Range("B5").Value = "2023-05-18 16:59:58"
$ws.Range("C5").Value = "Educacion"
$ws.Range("D5").Value = 10000

$ws.Range("A6").Value = "unab@unab.edu.co"
$ws.Range("B6").Value = "2023-05-18 17:00:14"
$ws.Range("C6").Value = "Alojamiento"
$ws.Range("D6").Value = 150000
